$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.817988034783468
$ws.Range("D2").Value = 0.314190933387863
$ws.Range("E2").Value = 0.3299995353571674
$ws.Range("F2").Value = 0.9014183446607476
$ws.Range("G2").Value = 0.3977407773541941
$ws.Range("H2").Value = 0.5587468790071313
$ws.Range("J2").Value = 0.4508836976495445
$ws.Range("L2").Value = 0.8676421976201993
$ws.Range("O2").Value = 1.858221274558872

$ws.Range("B3").Value = 1.682167639189515
$ws.Range("D3").Value = 0.3173200079630973
$ws.Range("E3").Value = 0.3261452349320777
$ws.Range("F3").Value = 0.9174211771114962
$ws.Range("G3").Value = 0.4009651538374612
$ws.Range("H3").Value = 0.5651110815936562
$ws.Range("J3").Value = 0.4367186947141164
$ws.Range("L3").Value = 0.7629942786872732
$ws.Range("O3").Value = 1.878068561783067

$ws.Range("B4").Value = 1.598625613031629
$ws.Range("D4").Value = 0.3194038750344177
$ws.Range("E4").Value = 0.3239447747195001
$ws.Range("F4").Value = 0.9280636962225586
$ws.Range("G4").Value = 0.4034045781797957
$ws.Range("H4").Value = 0.569395707108761
$ws.Range("J4").Value = 0.4282169067764983
$ws.Range("L4").Value = 0.6984489051218077
$ws.Range("O4").Value = 1.892007701312295

$ws.Range("B5").Value = 1.56454641203635
$ws.Range("D5").Value = 0.3202940530116862
$ws.Range("E5").Value = 0.323090093487572
$ws.Range("F5").Value = 0.9326055622688401
$ws.Range("G5").Value = 0.4045138315167591
$ws.Range("H5").Value = 0.5712364288630383
$ws.Range("J5").Value = 0.4248020524978671
$ws.Range("L5").Value = 0.6720747383461969
$ws.Range("O5").Value = 1.898127699385157

$ws.Range("B6").Value = 1.558885523820265
$ws.Range("D6").Value = 0.3204443446752592
$ws.Range("E6").Value = 0.3229507201525266
$ws.Range("F6").Value = 0.9333721020665529
$ws.Range("G6").Value = 0.4047049659075199
$ws.Range("H6").Value = 0.5715477963576561
$ws.Range("J6").Value = 0.4242380351804798
$ws.Range("L6").Value = 0.6676910692891624
$ws.Range("O6").Value = 1.899170444659504

$ws.Range("B7").Value = 1.59816614863297
$ws.Range("D7").Value = 0.3194157142090432
$ws.Range("E7").Value = 0.3239330776909028
$ws.Range("F7").Value = 0.9281241201257728
$ws.Range("G7").Value = 0.4034190721970958
$ws.Range("H7").Value = 0.5694201484056833
$ws.Range("J7").Value = 0.4281706509653844
$ws.Range("L7").Value = 0.6980935005702804
$ws.Range("O7").Value = 1.892088458913889

$ws.Range("B8").Value = 1.771188978019666
$ws.Range("D8").Value = 0.315236148318121
$ws.Range("E8").Value = 0.3286362505170359
$ws.Range("F8").Value = 0.9067663635454721
$ws.Range("G8").Value = 0.3987568947961293
$ws.Range("H8").Value = 0.5608629587050871
$ws.Range("J8").Value = 0.4459593106647048
$ws.Range("L8").Value = 0.8316210582425754
$ws.Range("O8").Value = 1.864700173703184

$ws.Range("B9").Value = 2.1092440944131
$ws.Range("D9").Value = 0.3083259704535237
$ws.Range("E9").Value = 0.3391667999297141
$ws.Range("F9").Value = 0.8713824610918621
$ws.Range("G9").Value = 0.3932803373824072
$ws.Range("H9").Value = 0.547077690442805
$ws.Range("J9").Value = 0.4823748098715015
$ws.Range("L9").Value = 1.091092868614851
$ws.Range("O9").Value = 1.824948853152506

$ws.Range("B10").Value = 2.356789285486229
$ws.Range("D10").Value = 0.3040273635303237
$ws.Range("E10").Value = 0.3476885719018128
$ws.Range("F10").Value = 0.849372320528893
$ws.Range("G10").Value = 0.3915180854046127
$ws.Range("H10").Value = 0.5387815414050436
$ws.Range("J10").Value = 0.5100404575821074
$ws.Range("L10").Value = 1.280211325162043
$ws.Range("O10").Value = 1.804321010042287

$ws.Range("B11").Value = 2.46921332655586
$ws.Range("D11").Value = 0.3022396465021941
$ws.Range("E11").Value = 0.351733392340094
$ws.Range("F11").Value = 0.8402301951073881
$ws.Range("G11").Value = 0.3912131090376505
$ws.Range("H11").Value = 0.5354064427724126
$ws.Range("J11").Value = 0.5228197217598165
$ws.Range("L11").Value = 1.365904686705903
$ws.Range("O11").Value = 1.796814156620968

$ws.Range("B12").Value = 2.511757229355283
$ws.Range("D12").Value = 0.3015867133084242
$ws.Range("E12").Value = 0.3532890322830013
$ws.Range("F12").Value = 0.8368939147344534
$ws.Range("G12").Value = 0.3911694920029589
$ws.Range("H12").Value = 0.5341858411503324
$ws.Range("J12").Value = 0.5276863748286758
$ws.Range("L12").Value = 1.398304589356087
$ws.Range("O12").Value = 1.794242571202432

$ws.Range("B13").Value = 2.502595943346194
$ws.Range("D13").Value = 0.3017262664578055
$ws.Range("E13").Value = 0.3529529357209213
$ws.Range("F13").Value = 0.8376068467058104
$ws.Range("G13").Value = 0.3911756821184014
$ws.Range("H13").Value = 0.5344461614426308
$ws.Range("J13").Value = 0.5266370418448787
$ws.Range("L13").Value = 1.391328956107827
$ws.Range("O13").Value = 1.794784331873387

$ws.Range("B14").Value = 2.472714024166692
$ws.Range("D14").Value = 0.3021854480480286
$ws.Range("E14").Value = 0.3518608970301003
$ws.Range("F14").Value = 0.8399531964674551
$ws.Range("G14").Value = 0.3912080776776321
$ws.Range("H14").Value = 0.5353048705772494
$ws.Range("J14").Value = 0.5232195573307763
$ws.Range("L14").Value = 1.36857126456141
$ws.Range("O14").Value = 1.796597150663473

$ws.Range("B15").Value = 2.454406698280536
$ws.Range("D15").Value = 0.3024698377997481
$ws.Range("E15").Value = 0.3511951044283137
$ws.Range("F15").Value = 0.8414067802740988
$ws.Range("G15").Value = 0.391237293875804
$ws.Range("H15").Value = 0.5358383431823768
$ws.Range("J15").Value = 0.5211298067114853
$ws.Range("L15").Value = 1.354624917166404
$ws.Range("O15").Value = 1.797742895059059

$ws.Range("B16").Value = 2.3494382080425
$ws.Range("D16").Value = 0.3041475631666728
$ws.Range("E16").Value = 0.3474275981592498
$ws.Range("F16").Value = 0.8499873358197476
$ws.Range("G16").Value = 0.3915480492463672
$ws.Range("H16").Value = 0.5390101483530714
$ws.Range("J16").Value = 0.5092091658007405
$ws.Range("L16").Value = 1.274604117340743
$ws.Range("O16").Value = 1.804849471099885

$ws.Range("B17").Value = 2.284994513212268
$ws.Range("D17").Value = 0.3052196952787938
$ws.Range("E17").Value = 0.345159271536815
$ws.Range("F17").Value = 0.8554745023286401
$ws.Range("G17").Value = 0.3918662197345242
$ws.Range("H17").Value = 0.5410581910822287
$ws.Range("J17").Value = 0.5019455894772733
$ws.Range("L17").Value = 1.225426215474329
$ws.Range("O17").Value = 1.809690694733547

$ws.Range("B18").Value = 2.247910858152295
$ws.Range("D18").Value = 0.3058521508401455
$ws.Range("E18").Value = 0.3438704378414883
$ws.Range("F18").Value = 0.8587124727071824
$ws.Range("G18").Value = 0.3920959462939209
$ws.Range("H18").Value = 0.5422737067081869
$ws.Range("J18").Value = 0.4977860625912172
$ws.Range("L18").Value = 1.197108716461059
$ws.Range("O18").Value = 1.812651825375951

$ws.Range("B19").Value = 2.235352051136317
$ws.Range("D19").Value = 0.306069004574411
$ws.Range("E19").Value = 0.3434367900510296
$ws.Range("F19").Value = 0.8598228450599592
$ws.Range("G19").Value = 0.3921817388129796
$ws.Range("H19").Value = 0.5426917025490781
$ws.Range("J19").Value = 0.4963808743217868
$ws.Range("L19").Value = 1.187515513002666
$ws.Range("O19").Value = 1.813684705385384

$ws.Range("B20").Value = 2.291856464192222
$ws.Range("D20").Value = 0.3051039310292509
$ws.Range("E20").Value = 0.3453991003179482
$ws.Range("F20").Value = 0.8548819045170717
$ws.Range("G20").Value = 0.3918275109944744
$ws.Range("H20").Value = 0.5408362877742974
$ws.Range("J20").Value = 0.5027169201102595
$ws.Range("L20").Value = 1.23066457782619
$ws.Range("O20").Value = 1.809157053360195

$ws.Range("B21").Value = 2.481491860947699
$ws.Range("D21").Value = 0.302049923609566
$ws.Range("E21").Value = 0.3521810070101949
$ws.Range("F21").Value = 0.8392606028936953
$ws.Range("G21").Value = 0.3911966081987401
$ws.Range("H21").Value = 0.5350510860883304
$ws.Range("J21").Value = 0.5242226148667442
$ws.Range("L21").Value = 1.375257129640147
$ws.Range("O21").Value = 1.796057314760191

$ws.Range("B22").Value = 2.605261729085385
$ws.Range("D22").Value = 0.3001940177069926
$ws.Range("E22").Value = 0.3567528927414685
$ws.Range("F22").Value = 0.8297837509997805
$ws.Range("G22").Value = 0.3912034015905306
$ws.Range("H22").Value = 0.5316051853437926
$ws.Range("J22").Value = 0.538437458654073
$ws.Range("L22").Value = 1.469462660427041
$ws.Range("O22").Value = 1.789076594363166

$ws.Range("B23").Value = 2.539219415189336
$ws.Range("D23").Value = 0.3011717612929701
$ws.Range("E23").Value = 0.3543001020841814
$ws.Range("F23").Value = 0.8347745363259236
$ws.Range("G23").Value = 0.3911612783350762
$ws.Range("H23").Value = 0.5334136299353673
$ws.Range("J23").Value = 0.5308362777044238
$ws.Range("L23").Value = 1.419210877329135
$ws.Range("O23").Value = 1.792657303619194

$ws.Range("B24").Value = 2.28875428226587
$ws.Range("D24").Value = 0.3051562179735896
$ws.Range("E24").Value = 0.3452906261803577
$ws.Range("F24").Value = 0.8551495584679998
$ws.Range("G24").Value = 0.3918448654511195
$ws.Range("H24").Value = 0.5409364916850166
$ws.Range("J24").Value = 0.5023681504085005
$ws.Range("L24").Value = 1.228296452601739
$ws.Range("O24").Value = 1.809397758729261

$ws.Range("B25").Value = 2.017931275546175
$ws.Range("D25").Value = 0.3100582761664938
$ws.Range("E25").Value = 0.336179296305481
$ws.Range("F25").Value = 0.8802564619558169
$ws.Range("G25").Value = 0.3943666744872232
$ws.Range("H25").Value = 0.5504856911070277
$ws.Range("J25").Value = 0.4723618386709063
$ws.Range("L25").Value = 1.021160286872544
$ws.Range("O25").Value = 1.83420129516611
